$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.647.50'
$ws.Range('E2').Value = '  +2.11%  '

$ws.Range('D3').Value = '3.920.04'
$ws.Range('E3').Value = '  +1.28%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '480.80'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.36%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '144.70'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.37%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.621'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.94%  '

$ws.Range('E8').Value = '  -0.17%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.722'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.83%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.168'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +8.47%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000351'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +11.99%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '42.55'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.12%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.48'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.53%  '

$ws.Range('D14').Value = '4.554.06'
$ws.Range('E14').Value = '  +0.87%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.942.21'
$ws.Range('E15').Value = '  +1.07%  '

$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.58'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.54%  '

$ws.Range('E17').Value = '  -0.32%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '19.66'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.10%  '

$ws.Range('E19').Value = '  -2.76%  '

$ws.Range('D20').Value = '68.698.25'
$ws.Range('E20').Value = '  +1.90%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '432.60'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.04%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '14.60'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.74%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.34'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.26%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '87.35'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.14%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.66'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +17.19%  '

$ws.Range('E26').Value = '  -0.99%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '38.11'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.54%  '

$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.17'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.43%  '

$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.82'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +6.04%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '706.20'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.29%  '

$ws.Range('E31').Value = '  -3.25%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '13.25'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.29%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.85'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.07%  '

$ws.Range('D34').Value = '0.0₃0912'
$ws.Range('E34').Value = '  +32.19%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '41.25'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -5.94%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '58.64'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.29%  '

$ws.Range('E37').Value = '  -7.53%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.63'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.90%  '

$ws.Range('E39').Value = '  -0.07%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0471'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.39%  '

$ws.Range('E41').Value = '  +9.30%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.74'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +6.75%  '

$ws.Range('E43').Value = '  +3.06%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.339'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.18%  '

$ws.Range('E46').Value = '  -0.02%  '

$ws.Range('E47').Value = '  -1.34%  '

$ws.Range('E48').Value = '  -0.80%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '147.59'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.22%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.14'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.17%  '

$ws.Range('E51').Value = '  -1.83%  '
